$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename shared-string header used by column O (effect1.invokeNum -> effect1.invokeTime)
$ws.Range("O1").Value = "effect1.invokeTime"

# --- Column L (protectAbility) was blank on most data rows; fill explicit 0s
$ws.Range("L2").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("L21").Value = 0

# --- Row 21 (hunyuan_sword_icon_big) attackThreaten changed from 0 to 10
$ws.Range("K21").Value = 10

# --- Selection / scroll position moved from K21 to M13, viewport scrolled so column F is leftmost
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$null = $ws.Range("M13").Select()
